$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in the "Fecha"/"Variedad"/"Volumen"/price/unit columns cycles between
# rows 2, 4 and 6: new row2 = old row6, new row4 = old row2, new row6 = old row4.
# Read the original ("before") values first so the rotation is computed correctly.
$cols = @("D", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

$row2 = @{}
$row4 = @{}
$row6 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col" + "2").Value()
    $row4[$col] = $ws.Range("$col" + "4").Value()
    $row6[$col] = $ws.Range("$col" + "6").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value = $row6[$col]
    $ws.Range("$col" + "4").Value = $row2[$col]
    $ws.Range("$col" + "6").Value = $row4[$col]
}
